$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B4"  = 0
    "B5"  = 3
    "B7"  = 12
    "B8"  = 32
    "B10" = 0
    "B11" = 27
    "B12" = 1
    "B13" = 206
    "B14" = 31
    "B15" = 2
    "B20" = 12
    "B21" = 2
    "B22" = 11
    "B27" = 53
    "B29" = 42
    "B30" = 5
    "B31" = 119
    "B32" = 109
    "B33" = 435
    "B35" = 28
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
